$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 319 (shifts rows 319..376 down to 320..377)
$ws.Rows.Item(319).Insert()

# Fill in the new row 319 with its values.
# Columns A,B,C,E,F,G,I,N,O,Q,R keep the same values the (old) row 319 had;
# columns D,H,J,K,L,M,P get new values per the target data.
$ws.Cells.Item(319, 1).Value = 7
$ws.Cells.Item(319, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(319, 3).Value = "Ñuble"
$ws.Cells.Item(319, 4).Value = "2023-09-11"
$ws.Cells.Item(319, 5).Value = 16
$ws.Cells.Item(319, 6).Value = 100112045
$ws.Cells.Item(319, 7).Value = "Zapallo"
$ws.Cells.Item(319, 8).Value = "Paine"
$ws.Cells.Item(319, 9).Value = "1a (guarda)"
$ws.Cells.Item(319, 10).Value = 270
$ws.Cells.Item(319, 11).Value = 350
$ws.Cells.Item(319, 12).Value = 400
$ws.Cells.Item(319, 13).Value = 378
$ws.Cells.Item(319, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(319, 15).Value = "Región del Maule"
$ws.Cells.Item(319, 16).Value = 378
$ws.Cells.Item(319, 17).Value = 1
$ws.Cells.Item(319, 18).Value = "Hortaliza"
